$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 29.083334
$ws.Range("I11").Value = 29.083334
$ws.Range("K11").Value = 29.083334
$ws.Range("M11").Value = 110.916666
$ws.Range("H62").Value = 201240.2
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 201240.2
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H76").Value = 2598.8823
$ws.Range("I76").Value = 2463.3572
$ws.Range("K76").Value = 2463.3572
$ws.Range("M76").Value = -2148.3572
$ws.Range("H79").Value = 2598.8823
$ws.Range("I79").Value = 2463.3572
$ws.Range("K79").Value = 2463.3572
$ws.Range("M79").Value = -1371.3572
$ws.Range("H103").Value = 4751.25
$ws.Range("J103").Value = 4751.25
$ws.Range("L103").Value = 14253.75
$ws.Range("N103").Value = -15425.75
$ws.Range("H135").Value = 5881.7393
$ws.Range("I135").Value = 2460.353
$ws.Range("K135").Value = 22143.177
$ws.Range("M135").Value = -19608.177

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5470.4194
$ws.Range("I61").Value = 4578.5835
$ws.Range("K61").Value = 4578.5835
$ws.Range("M61").Value = -4366.5835
$ws.Range("H74").Value = 1324.2632
$ws.Range("I74").Value = 1105.8182
$ws.Range("K74").Value = 1105.8182
$ws.Range("M74").Value = -231.8181999999999
$ws.Range("H77").Value = 1324.2632
$ws.Range("I77").Value = 1105.8182
$ws.Range("K77").Value = 5529.090999999999
$ws.Range("M77").Value = -1161.090999999999
$ws.Range("H132").Value = 7458.933
$ws.Range("I132").Value = 4231.6665
$ws.Range("J132").Value = 12299.833
$ws.Range("K132").Value = 12694.9995
$ws.Range("L132").Value = 36899.499
$ws.Range("M132").Value = -10164.9995
$ws.Range("N132").Value = -41959.499
$ws.Range("H136").Value = 5470.4194
$ws.Range("I136").Value = 4578.5835
$ws.Range("K136").Value = 13735.7505
$ws.Range("M136").Value = -11185.7505

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1430.1212
$ws.Range("I20").Value = 1247.6
$ws.Range("K20").Value = 1247.6
$ws.Range("M20").Value = -1000.6
$ws.Range("H22").Value = 450.08334
$ws.Range("I22").Value = 400.1111
$ws.Range("K22").Value = 400.1111
$ws.Range("M22").Value = -227.1111
$ws.Range("H99").Value = 2960
$ws.Range("I99").Value = 3107.5386
$ws.Range("K99").Value = 3107.5386
$ws.Range("M99").Value = -1609.5386
$ws.Range("H105").Value = 3116.4443
$ws.Range("I105").Value = 2855.6667
$ws.Range("K105").Value = 2855.6667
$ws.Range("M105").Value = -1108.6667
$ws.Range("H134").Value = 2968.9512
$ws.Range("I134").Value = 2398.1052
$ws.Range("K134").Value = 7194.3156
$ws.Range("M134").Value = -4659.3156

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2591.3684
$ws.Range("I31").Value = 1115.65
$ws.Range("K31").Value = 1115.65
$ws.Range("M31").Value = -820.6500000000001
$ws.Range("H34").Value = 2591.3684
$ws.Range("I34").Value = 1115.65
$ws.Range("K34").Value = 1115.65
$ws.Range("M34").Value = -913.6500000000001
$ws.Range("H58").Value = 1984.1
$ws.Range("J58").Value = 3541.6667
$ws.Range("L58").Value = 3541.6667
$ws.Range("N58").Value = -3947.6667
$ws.Range("H99").Value = 2928.7297
$ws.Range("I99").Value = 2928.7297
$ws.Range("K99").Value = 2928.7297
$ws.Range("M99").Value = -1430.7297
$ws.Range("H107").Value = 1458.2
$ws.Range("I107").Value = 989.1818
$ws.Range("J107").Value = 2031.4445
$ws.Range("K107").Value = 989.1818
$ws.Range("L107").Value = 2031.4445
$ws.Range("M107").Value = 930.8182
$ws.Range("N107").Value = -5871.4445
$ws.Range("H122").Value = 2724.3125
$ws.Range("I122").Value = 1915.2
$ws.Range("J122").Value = 5614
$ws.Range("K122").Value = 5745.6
$ws.Range("L122").Value = 16842
$ws.Range("M122").Value = -3295.6
$ws.Range("N122").Value = -21742
$ws.Range("H126").Value = 2928.7297
$ws.Range("I126").Value = 2928.7297
$ws.Range("K126").Value = 8786.1891
$ws.Range("M126").Value = -6316.1891
$ws.Range("H132").Value = 1298.7273
$ws.Range("I132").Value = 1298.7273
$ws.Range("K132").Value = 3896.1819
$ws.Range("M132").Value = -1366.1819
$ws.Range("H136").Value = 1984.1
$ws.Range("J136").Value = 3541.6667
$ws.Range("L136").Value = 10625.0001
$ws.Range("N136").Value = -15725.0001
$ws.Range("H141").Value = 86125.55499999999
$ws.Range("I141").Value = 84644
$ws.Range("J141").Value = 87310.8
$ws.Range("K141").Value = 84644
$ws.Range("L141").Value = 87310.8
$ws.Range("M141").Value = -79464
$ws.Range("N141").Value = -97670.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 426957.3
$ws.Range("I56").Value = 426957.3
$ws.Range("K56").Value = 426957.3
$ws.Range("M56").Value = -426427.3
$ws.Range("H97").Value = 1505
$ws.Range("I97").Value = 2819.8
$ws.Range("K97").Value = 8459.400000000001
$ws.Range("M97").Value = -7963.400000000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 92351.75
$ws.Range("J69").Value = 92351.75
$ws.Range("L69").Value = 92351.75
$ws.Range("N69").Value = -93849.75
$ws.Range("H72").Value = 92351.75
$ws.Range("J72").Value = 92351.75
$ws.Range("L72").Value = 277055.25
$ws.Range("N72").Value = -284543.25
$ws.Range("H80").Value = 4698.2
$ws.Range("I80").Value = 4204.3105
$ws.Range("K80").Value = 4204.3105
$ws.Range("M80").Value = -3206.3105
$ws.Range("H83").Value = 4698.2
$ws.Range("I83").Value = 4204.3105
$ws.Range("K83").Value = 21021.5525
$ws.Range("M83").Value = -16029.5525
$ws.Range("H102").Value = 128405.5
$ws.Range("I102").Value = 4541.3335
$ws.Range("K102").Value = 4541.3335
$ws.Range("M102").Value = -2919.3335
$ws.Range("H105").Value = 82999.664
$ws.Range("J105").Value = 82999.664
$ws.Range("L105").Value = 82999.664
$ws.Range("N105").Value = -89987.664
$ws.Range("H126").Value = 4446.037
$ws.Range("I126").Value = 3785.225
$ws.Range("J126").Value = 6334.0713
$ws.Range("K126").Value = 11355.675
$ws.Range("L126").Value = 19002.2139
$ws.Range("M126").Value = -8885.674999999999
$ws.Range("N126").Value = -23942.2139
$ws.Range("H132").Value = 6976
$ws.Range("I132").Value = 7158.722
$ws.Range("K132").Value = 21476.166
$ws.Range("M132").Value = -18946.166

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4184.3125
$ws.Range("I7").Value = 3412.5833
$ws.Range("J7").Value = 6499.5
$ws.Range("K7").Value = 3412.5833
$ws.Range("L7").Value = 6499.5
$ws.Range("M7").Value = -3300.5833
$ws.Range("N7").Value = -6723.5
$ws.Range("H40").Value = 6399.769
$ws.Range("I40").Value = 5417.909
$ws.Range("J40").Value = 11800
$ws.Range("K40").Value = 5417.909
$ws.Range("L40").Value = 11800
$ws.Range("M40").Value = -5281.909
$ws.Range("N40").Value = -12072
$ws.Range("H122").Value = 13648.976
$ws.Range("J122").Value = 13272.385
$ws.Range("L122").Value = 39817.155
$ws.Range("N122").Value = -44717.155
$ws.Range("H126").Value = 4184.3125
$ws.Range("I126").Value = 3412.5833
$ws.Range("J126").Value = 6499.5
$ws.Range("K126").Value = 10237.7499
$ws.Range("L126").Value = 19498.5
$ws.Range("M126").Value = -7767.749899999999
$ws.Range("N126").Value = -24438.5
$ws.Range("H136").Value = 3741.081
$ws.Range("I136").Value = 3758.2122
$ws.Range("K136").Value = 11274.6366
$ws.Range("M136").Value = -8724.6366

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12292.412
$ws.Range("I122").Value = 13717.889
$ws.Range("J122").Value = 10688.75
$ws.Range("K122").Value = 41153.667
$ws.Range("L122").Value = 32066.25
$ws.Range("M122").Value = -38703.667
$ws.Range("N122").Value = -36966.25
$ws.Range("H126").Value = 9660.267
$ws.Range("J126").Value = 16399.2
$ws.Range("L126").Value = 49197.60000000001
$ws.Range("N126").Value = -54137.60000000001
$ws.Range("H132").Value = 1900.7778
$ws.Range("I132").Value = 1800.75
$ws.Range("K132").Value = 5402.25
$ws.Range("M132").Value = -2872.25
$ws.Range("H136").Value = 3978.1052
$ws.Range("I136").Value = 2634.261
$ws.Range("K136").Value = 7902.782999999999
$ws.Range("M136").Value = -5352.782999999999
